$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 - Bitcoin
$ws.Range("D2").Value = "72.171.48"

# Row 3 - Ethereum
$ws.Range("D3").Value = "4.040.93"
$ws.Range("E3").Value = "  +0.03%  "

# Row 4 - TetherUSD
$ws.Range("E4").Value = "  -0.03%  "

# Row 5 - BNB
$ws.Range("D5").Value = "'539.21"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -0.01%  "

# Row 6 - Solana
$ws.Range("D6").Value = "'149.62"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -2.74%  "

# Row 7 - LidoStakedEther
$ws.Range("D7").Value = "4.034.59"
$ws.Range("E7").Value = "  +0.01%  "

# Row 8 - XRP
$ws.Range("D8").Value = "'0.696"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  +0.56%  "

# Row 9 - USDC
$ws.Range("E9").Value = "  -0.06%  "

# Row 10 - Cardano
$ws.Range("E10").Value = "  -1.30%  "

# Row 11 - Dogecoin
$ws.Range("E11").Value = "  -0.70%  "

# Row 12 - Avalanche
$ws.Range("D12").Value = "'53.58"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  +11.70%  "

# Row 13 - ShibaInu
$ws.Range("E13").Value = "  +0.51%  "

# Row 14 - Polkadot
$ws.Range("D14").Value = "'10.86"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  -0.11%  "

# Row 15 - WrappedliquidstakedEther2.0
$ws.Range("D15").Value = "4.683.54"
$ws.Range("E15").Value = "  -0.13%  "

# Row 16 - WrappedEther
$ws.Range("D16").Value = "4.038.03"
$ws.Range("E16").Value = "  +0.05%  "

# Row 17 - Uniswap
$ws.Range("D17").Value = "'14.32"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  -0.38%  "

# Row 18 - Chainlink
$ws.Range("D18").Value = "'20.63"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  -0.15%  "

# Row 19 - Polygon
$ws.Range("E19").Value = "  -0.83%  "

# Row 20 - TRON
$ws.Range("E20").Value = "  -0.92%  "

# Row 21 - WrappedBTC
$ws.Range("D21").Value = "72.137.56"
$ws.Range("E21").Value = "  +0.39%  "

# Row 22 - BitcoinCash
$ws.Range("D22").Value = "'440.11"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +1.18%  "

# Row 23 - Litecoin
$ws.Range("D23").Value = "'97.65"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  -1.45%  "

# Row 24 - ImmutableX
$ws.Range("E24").Value = "  -2.75%  "

# Row 25 - PancakeSwap
$ws.Range("D25").Value = "'4.27"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  -0.07%  "

# Row 26 - InternetComputer(DFINITY)
$ws.Range("D26").Value = "'14.62"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  -0.72%  "

# Row 27 - Toncoin
$ws.Range("D27").Value = "'4.45"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  +23.34%  "

# Row 28 - RenderToken
$ws.Range("D28").Value = "'11.22"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  -0.78%  "

# Row 29 - Filecoin
$ws.Range("E29").Value = "  -2.93%  "

# Row 30 - LEO
$ws.Range("E30").Value = "  +1.99%  "

# Row 31 - EthereumClassic
$ws.Range("D31").Value = "'37.19"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  +0.32%  "

# Row 32 - NEARProtocol
$ws.Range("D32").Value = "'8.32"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  +20.69%  "

# Row 33 - Hedera
$ws.Range("E33").Value = "  +1.95%  "

# Row 34 - Cosmos
$ws.Range("E34").Value = "  -0.84%  "

# Row 35 - InjectiveProtocol
$ws.Range("D35").Value = "'49.26"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  +14.23%  "

# Row 36 - Bittensor
$ws.Range("D36").Value = "'684.97"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  -0.58%  "

# Row 37 - OKB
$ws.Range("D37").Value = "'67.19"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  +0.27%  "

# Row 38 - TheGraph
$ws.Range("D38").Value = "'0.457"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  +4.75%  "

# Row 39 - PEPE
$ws.Range("D39").Value = "0.0₃0911"
$ws.Range("E39").Value = "  +8.16%  "

# Row 40 - Kaspa (was THORChain)
$ws.Range("B40").Value = "Kaspa"
$ws.Range("C40").Value = "https://coinranking.com/coin/V8GxkwWow+kaspa-kas"
$ws.Range("D40").Value = "'0.148"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  -6.32%  "

# Row 41 - THORChain (was Kaspa)
$ws.Range("B41").Value = "THORChain"
$ws.Range("C41").Value = "https://coinranking.com/coin/ybmU-kKU+thorchain-rune"
$ws.Range("D41").Value = "'11.34"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  +18.61%  "

# Row 42 - WEMIXToken (was ThetaToken)
$ws.Range("B42").Value = "WEMIXToken"
$ws.Range("C42").Value = "https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix"
$ws.Range("E42").Value = "  +3.38%  "

# Row 43 - ThetaToken (was WEMIXToken)
$ws.Range("B43").Value = "ThetaToken"
$ws.Range("C43").Value = "https://coinranking.com/coin/B42IRxNtoYmwK+thetatoken-theta"
$ws.Range("E43").Value = "  -1.55%  "

# Row 44 - Dai
$ws.Range("E44").Value = "  +0.10%  "

# Row 45 - FirstDigitalUSD
$ws.Range("D45").Value = "'1.00"
$ws.Range("D45").Style = "Normal"

# Row 46 - VeChain
$ws.Range("D46").Value = "'0.0491"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  -0.74%  "

# Row 47 - Stellar
$ws.Range("D47").Value = "'0.151"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -0.90%  "

# Row 48 - Fetch.AI
$ws.Range("D48").Value = "'2.63"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  -3.29%  "

# Row 49 - Stacks
$ws.Range("E49").Value = "  +2.22%  "

# Row 50 - ApeXProtocol (was FLOKI)
$ws.Range("B50").Value = "ApeXProtocol"
$ws.Range("C50").Value = "https://coinranking.com/coin/ze0N2Rcyu+apexprotocol-apex"
$ws.Range("D50").Value = "'3.35"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -1.69%  "

# Row 51 - FLOKI (was ApeXProtocol)
$ws.Range("B51").Value = "FLOKI"
$ws.Range("C51").Value = "https://coinranking.com/coin/fmHk13Rqw+floki-floki"
$ws.Range("D51").Value = "'0.000283"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +5.00%  "
